$d = $word.ActiveDocument

# Position a collapsed range at the very end of the document body
# (immediately after the "Cccccccc" paragraph) and inject the two
# new paragraphs as raw WordprocessingML, so their run formatting is
# taken verbatim from the markup instead of being inherited from the
# surrounding text.
$endRange = $d.Content
$endRange.Collapse(0)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>D</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ddddd</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Eeeeeee</w:t></w:r></w:p>'

$null = $endRange.InsertXML($newParagraphsXml)
